$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "Hola soy marc" + " holaaaaaa"
# We rebuild it (splitting "marc" and "holaaaaaa" out into their own runs
# wrapped in spell-check proofErr markers, matching what Word's spell
# checker leaves behind), and append two further paragraphs: an empty
# paragraph, and a new "Holaaa soy pablo" paragraph (also spell-marked)
# that inherits the trailing _GoBack bookmark.

$p1 = $d.Paragraphs(1).Range
$startPos = $p1.Start
$endPos = $p1.End - 1   # End() includes the paragraph mark; stop right before it

$target = $d.Range($startPos, $endPos)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:r><w:t xml:space="preserve">Hola soy </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>marc</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>holaaaaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
'</w:p>' +
'<w:p/>' +
'<w:p>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Holaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> soy pablo</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newXml)
